$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the standalone "Meta description" paragraph that sits
#    right after the Heading1 title at the top of the document.
#    Locate it by content (instead of a hard-coded paragraph index)
#    so the script is resilient to minor structural differences.
# -----------------------------------------------------------------
$metaFind = $d.Content
$metaFound = $metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "", 0)
if ($metaFound) {
    $metaPara = $metaFind.Paragraphs(1).Range
    $metaPara.Delete()
}

# -----------------------------------------------------------------
# 2) At the end of the document, insert a new bold "title" paragraph
#    right before the final (italic) paragraph, and rewrite that
#    final paragraph's text (dropping the old DALLE image prompt,
#    replacing it with the meta-description copy) while keeping its
#    italic formatting.
# -----------------------------------------------------------------
$promptFind = $d.Content
$promptFound = $promptFind.Find.Execute("Prompt: Create a feature image", $true, $false, $false, $false, `
                                         $false, $true, 1, $false, "", 0)
if ($promptFound) {
    $targetPara = $promptFind.Paragraphs(1).Range
} else {
    $targetPara = $d.Paragraphs($d.Paragraphs.Count).Range
}

# Work against the paragraph's text only (exclude its trailing
# paragraph mark) so the document's final paragraph mark is never
# disturbed -- Word will not let that last mark be deleted/merged
# away, and doing so would leave a stray empty trailing paragraph.
$full = $targetPara.Paragraphs(1).Range
$r = $d.Range($full.Start, $full.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Bad Wolf for Free - A Fairytale Themed Slot Game</w:t></w:r></w:p>
          <w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Big Bad Wolf, an online slot game developed by Quickspin and inspired by The Three Little Pigs story. Play this fairy tale-themed game for free now!</w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
